# Uno_v1.0 PnP workbook — fix rotation values and refresh last-used selection
# (per commit: "Fixed unconnected symbol on schematics... Fixed BoM and PnP
# files for JLCPB assembly")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Rotation ("Layer"/"Rotation" columns) corrections for re-oriented parts
$ws.Range("E25").Value = 270
$ws.Range("E26").Value = 270
$ws.Range("E51").Value = 270
$ws.Range("E52").Value = 90
$ws.Range("E54").Value = 180

# Leave the cursor where the editor last left it when saving
$ws.Range("E53").Select()
